# Apply updated cryptocurrency price/volume data to Sheet1
# (values mirror the upstream coinranking.com scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.280.90'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '1.576.30'
$ws.Range('E3').Value = '  -1.01%  '
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').Value = '''208.09'
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('D6').Value = '''0.491'
$ws.Range('E6').Value = '  -2.09%  '
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').Value = '''22.28'
$ws.Range('E9').Value = '  -1.47%  '
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('D12').Value = '1.800.44'
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('D13').Value = '1.577.45'
$ws.Range('E13').Value = '  -1.06%  '
$ws.Range('E14').Value = '  -1.40%  '
$ws.Range('D15').Value = '''0.521'
$ws.Range('E15').Value = '  -1.63%  '
$ws.Range('D16').Value = '''62.62'
$ws.Range('E16').Value = '  -1.04%  '
$ws.Range('D17').Value = '27.272.22'
$ws.Range('E17').Value = '  -1.61%  '
$ws.Range('D18').Value = '''215.89'
$ws.Range('D19').Value = '''7.32'
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('E20').Value = '  -1.10%  '
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('E23').Value = '  -3.29%  '
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('D25').Value = '''151.39'
$ws.Range('E25').Value = '  -1.60%  '
$ws.Range('E26').Value = '  -5.36%  '
$ws.Range('E27').Value = '  -1.07%  '
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('E30').Value = '  -1.68%  '
$ws.Range('E31').Value = '  -2.32%  '
$ws.Range('E32').Value = '  -1.13%  '
$ws.Range('D33').Value = '1.406.53'
$ws.Range('E33').Value = '  +1.77%  '
$ws.Range('E34').Value = '  -1.66%  '
$ws.Range('E35').Value = '  +1.40%  '
$ws.Range('E36').Value = '  -2.29%  '
$ws.Range('D37').Value = '''0.940'
$ws.Range('E37').Value = '  -3.04%  '
$ws.Range('D38').Value = '''0.0165'
$ws.Range('E38').Value = '  -1.93%  '
$ws.Range('D39').Value = '''0.821'
$ws.Range('E39').Value = '  -0.72%  '
$ws.Range('E40').Value = '  -2.67%  '
$ws.Range('E42').Value = '  +1.75%  '
$ws.Range('E43').Value = '  +3.39%  '
$ws.Range('E44').Value = '  +1.74%  '
$ws.Range('E45').Value = '  +0.61%  '
$ws.Range('D46').Value = '''63.83'
$ws.Range('E46').Value = '  -1.01%  '
$ws.Range('D47').Value = '1.712.72'
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('D48').Value = '''86.23'
$ws.Range('E48').Value = '  +0.17%  '
$ws.Range('D49').Value = '0.0₇0989'
$ws.Range('E49').Value = '  -1.74%  '
$ws.Range('D50').Value = '''0.0954'
$ws.Range('E50').Value = '  -1.49%  '
$ws.Range('E51').Value = '  -0.41%  '
